# Auto-generated edit script applying the crypto price-update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.552.41"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "1.880.86"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7271"
$ws.Range("E5").Value = "  +3.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.39"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9988"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07879"
$ws.Range("E8").Value = "  -3.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3089"
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.32"
$ws.Range("E10").Value = "  +9.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08215"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "1.868.07"
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.273"
$ws.Range("E13").Value = "  +2.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7271"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.71"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").Value = "29.541.34"
$ws.Range("E16").Value = "  +0.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.858"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007872"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.26"
$ws.Range("E19").Value = "  +2.91%  "
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.122.47"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9986"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9993"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.782"
$ws.Range("E24").Value = "  +5.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1494"
$ws.Range("E25").Value = "  +3.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.89"
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.004"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.958"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.362"
$ws.Range("E30").Value = "  -4.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.480"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.364"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.106"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05254"
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.200"
$ws.Range("E35").Value = "  +2.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7193"
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.003"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.670"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01862"
$ws.Range("E39").Value = "  +1.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.717"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").Value = "1.176.18"
$ws.Range("E41").Value = "  +3.85%  "
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.996"
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.05"
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4320"
$ws.Range("E45").Value = "  +1.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9989"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.58"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5349"
$ws.Range("E48").Value = "  -1.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.780"
$ws.Range("E49").Value = "  +0.99%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.266"
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.890"
$ws.Range("E51").Value = "  +5.46%  "

Write-Output "Applied 103 cell updates"
